# Update TPM-derived NATMI metrics for Nid1-Col13a1 ligand-receptor pairs (rows 2-10).
# Only numeric columns G..T are affected; columns A-F (identifiers/counts) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 28.743868
$ws.Range("H2").Value = 86.231604
$ws.Range("I2").Value = 0.0554303735704667
$ws.Range("J2").Value = 0.0554303735704667
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1837056666666667
$ws.Range("N2").Value = 0.5511170000000001
$ws.Range("O2").Value = 0.7269991860920679
$ws.Range("P2").Value = 0.7269991860920678
$ws.Range("Q2").Value = 5.280411433518668
$ws.Range("R2").Value = 47.52370290166801
$ws.Range("S2").Value = 0.04029783647050857
$ws.Range("T2").Value = 0.04029783647050856

# Row 3
$ws.Range("G3").Value = 28.743868
$ws.Range("H3").Value = 86.231604
$ws.Range("I3").Value = 0.0554303735704667
$ws.Range("J3").Value = 0.0554303735704667
$ws.Range("O3").Value = 0.2534828531892131
$ws.Range("P3").Value = 0.2534828531892131
$ws.Range("Q3").Value = 1.841121395714667
$ws.Range("R3").Value = 16.570092561432
$ws.Range("S3").Value = 0.01405064924598585
$ws.Range("T3").Value = 0.01405064924598585

# Row 4
$ws.Range("G4").Value = 28.743868
$ws.Range("H4").Value = 86.231604
$ws.Range("I4").Value = 0.0554303735704667
$ws.Range("J4").Value = 0.0554303735704667
$ws.Range("O4").Value = 0.01951796071871896
$ws.Range("P4").Value = 0.01951796071871896
$ws.Range("Q4").Value = 0.141764756976
$ws.Range("R4").Value = 1.275882812784
$ws.Range("S4").Value = 0.001081887853972286
$ws.Range("T4").Value = 0.001081887853972286

# Row 5
$ws.Range("I5").Value = 0.848161237947095
$ws.Range("J5").Value = 0.8481612379470951
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1837056666666667
$ws.Range("N5").Value = 0.5511170000000001
$ws.Range("O5").Value = 0.7269991860920679
$ws.Range("P5").Value = 0.7269991860920678
$ws.Range("Q5").Value = 80.79758460638278
$ws.Range("R5").Value = 727.178261457445
$ws.Range("S5").Value = 0.6166125296623788
$ws.Range("T5").Value = 0.6166125296623788

# Row 6
$ws.Range("I6").Value = 0.848161237947095
$ws.Range("J6").Value = 0.8481612379470951
$ws.Range("O6").Value = 0.2534828531892131
$ws.Range("P6").Value = 0.2534828531892131
$ws.Range("S6").Value = 0.2149943305593247
$ws.Range("T6").Value = 0.2149943305593247

# Row 7
$ws.Range("I7").Value = 0.848161237947095
$ws.Range("J7").Value = 0.8481612379470951
$ws.Range("O7").Value = 0.01951796071871896
$ws.Range("P7").Value = 0.01951796071871896
$ws.Range("S7").Value = 0.01655437772539144
$ws.Range("T7").Value = 0.01655437772539144

# Row 8
$ws.Range("G8").Value = 49.99334866666666
$ws.Range("I8").Value = 0.09640838848243828
$ws.Range("J8").Value = 0.09640838848243828
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1837056666666667
$ws.Range("N8").Value = 0.5511170000000001
$ws.Range("O8").Value = 0.7269991860920679
$ws.Range("P8").Value = 0.7269991860920678
$ws.Range("Q8").Value = 9.18406144570911
$ws.Range("R8").Value = 82.656553011382
$ws.Range("S8").Value = 0.07008881995918052
$ws.Range("T8").Value = 0.07008881995918051

# Row 9
$ws.Range("G9").Value = 49.99334866666666
$ws.Range("I9").Value = 0.09640838848243828
$ws.Range("J9").Value = 0.09640838848243828
$ws.Range("O9").Value = 0.2534828531892131
$ws.Range("P9").Value = 0.2534828531892131
$ws.Range("S9").Value = 0.02443787338390253
$ws.Range("T9").Value = 0.02443787338390253

# Row 10
$ws.Range("G10").Value = 49.99334866666666
$ws.Range("I10").Value = 0.09640838848243828
$ws.Range("J10").Value = 0.09640838848243828
$ws.Range("O10").Value = 0.01951796071871896
$ws.Range("P10").Value = 0.01951796071871896
$ws.Range("S10").Value = 0.001881695139355228
$ws.Range("T10").Value = 0.001881695139355228

Write-Host "Applied 85 cell updates across rows 2-10"
